$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing end time on row 5 (2014-02-18 entry, end time moved from 12:00 to 13:15)
$ws.Range("E5").Value = 0.55208333333333337

# Insert a new row above row 6, shifting the blank separator row and summary rows down
$ws.Rows.Item(6).Insert()

# Fill in the new row 6 with a new time-tracking entry (2014-02-18, 13:45 - 18:00)
$ws.Range("A6").Value = 2014
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 18
$ws.Range("D6").Value = 0.57291666666666663
$ws.Range("E6").Value = 0.75
$ws.Range("F6").Formula = "=(E6-D6)*24*60"

# Update the selected cell to match the final state
$ws.Range("L12").Select()
